$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "2003 budget" sheet: column A relabelling (post-vote corrections).
#    Rows 2-26 were "Income" -> now "In"
#    Rows 27-204 were "Expense" -> now "Out"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2003 budget")
$ws1.Range("A2:A26").Value = "In"
$ws1.Range("A27:A204").Value = "Out"

# ---------------------------------------------------------------------------
# 2. Defined-name corrections for the data-model worksheet connections
#    (Table2 -> Table21, Table3 -> Table31).
# ---------------------------------------------------------------------------
$n2 = $wb.Names.Item("_xlcn.WorksheetConnection_budget_2024.xlsxTable2")
$n2.Name = "_xlcn.WorksheetConnection_budget_2024.xlsxTable21"

$n3 = $wb.Names.Item("_xlcn.WorksheetConnection_budget_2024.xlsxTable3")
$n3.Name = "_xlcn.WorksheetConnection_budget_2024.xlsxTable31"

# ---------------------------------------------------------------------------
# 3. View / selection state: "2003 budget" becomes the active (tab-selected)
#    sheet, with a new selection on the newly re-labelled "Out" rows.
# ---------------------------------------------------------------------------
[void]$ws1.Activate()
[void]$ws1.Range("A27:A204").Select()
